# MSME Country Indicators - Kazakhstan Summary
# Inserts the "MSME definition" table (Micro/Small/Medium/Large, by number
# of employees / assets / turnover) above the existing "Sector Distribution
# Details" block, pushing that block (and everything below it) down by six
# rows, and repoints the DAMU source hyperlink at its new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 6 blank rows at row 19 (the former row 22, "Sector
# Distribution Details", ends up at row 28 - matching the 1 blank row that
# used to separate it from the "Source: KazStat, 2012" note at row 17).
$ws.Rows("19:24").Insert()

# --- New table header (row 19) ---
$ws.Range("B19").Value = "Number of employees"
$ws.Range("C19").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D19").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B19:D19").Font.Bold = $true

# --- Micro (row 20) ---
$ws.Range("A20").Value = "Micro"

# --- Small (row 21) ---
$ws.Range("A21").Value = "Small"
$ws.Range("B21").Value = "<50"
$ws.Range("C21").Value = "<60,000 conventional units"

# --- Medium (row 22) ---
$ws.Range("A22").Value = "Medium"
$ws.Range("B22").Value = "51-250"
$ws.Range("C22").Value = "<325,000 conventional units"

# --- Large (row 23) ---
$ws.Range("A23").Value = "Large"
$ws.Range("B23").Value = ">250"
$ws.Range("C23").Value = ">325,000 conventional units"

# The row insert above shifts cell contents but leaves the worksheet's
# hyperlink collection pointing at the old anchor (A36). Recreate it at the
# new location (A42) with the same external target.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A42"), "http://www.damu.kz/content/files/27_01_2011_Damu_ADB.pdf")
